$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The LR-pair table was regenerated against the new TPM expression matrix.
# The refreshed data now covers only 6 sending/target-cluster combinations
# (FAPs and MuSCs as sending clusters, against ECs/FAPs/MuSCs as target),
# instead of the original 9 (which also had an "ECs" sending-cluster block),
# so the trailing three rows are removed and the sheet shrinks to A1:T7.
$ws.Range("A8:T10").EntireRow.Delete()

# New per-row values (columns A:T) for the remaining data rows 2-7.
$newData = @{
    2 = @("FAPs", "Col2a1", "Tnfrsf10b", "ECs", 1, 0.3333333333333333, 0.2021233333333333, 0.60637, 0.9764223557676824, 0.9764223557676823, 3, 1, 5.857771333333333, 17.573314, 0.7247832978729415, 0.7247832978729416, 1.183992267797778, 10.65593041018, 0.7076946151301674, 0.7076946151301674)
    3 = @("FAPs", "Col2a1", "Tnfrsf10b", "FAPs", 1, 0.3333333333333333, 0.2021233333333333, 0.60637, 0.9764223557676824, 0.9764223557676823, 3, 1, 1.267983333333333, 3.803949999999999, 0.1568878485835839, 0.1568878485835839, 0.2562890179444444, 2.3066011615, 0.1531888027053065, 0.1531888027053065)
    4 = @("FAPs", "Col2a1", "Tnfrsf10b", "MuSCs", 1, 0.3333333333333333, 0.2021233333333333, 0.60637, 0.9764223557676824, 0.9764223557676823, 3, 1, 0.9563456666666665, 2.869037, 0.1183288535434745, 0.1183288535434745, 0.1932997739655555, 1.73969796569, 0.1155389379322085, 0.1155389379322084)
    5 = @("MuSCs", "Col2a1", "Tnfrsf10b", "ECs", 2, 0.6666666666666666, 0.004880666666666667, 0.014642, 0.02357764423231757, 0.02357764423231757, 3, 1, 5.857771333333333, 17.573314, 0.7247832978729415, 0.7247832978729416, 0.02858982928755556, 0.257308463588, 0.01708868274277407, 0.01708868274277407)
    6 = @("MuSCs", "Col2a1", "Tnfrsf10b", "FAPs", 2, 0.6666666666666666, 0.004880666666666667, 0.014642, 0.02357764423231757, 0.02357764423231757, 3, 1, 1.267983333333333, 3.803949999999999, 0.1568878485835839, 0.1568878485835839, 0.006188603988888888, 0.05569743589999999, 0.003699045878277451, 0.003699045878277451)
    7 = @("MuSCs", "Col2a1", "Tnfrsf10b", "MuSCs", 2, 0.6666666666666666, 0.004880666666666667, 0.014642, 0.02357764423231757, 0.02357764423231757, 3, 1, 0.9563456666666665, 2.869037, 0.1183288535434745, 0.1183288535434745, 0.004667604417111111, 0.042008439754, 0.002789915611266053, 0.002789915611266053)
}

foreach ($r in $newData.Keys) {
    $values = $newData[$r]
    for ($col = 1; $col -le $values.Length; $col++) {
        $ws.Cells.Item($r, $col).Value = $values[$col - 1]
    }
}
